$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_dede = $wb.Worksheets.Item("de-de")

# --- Text change: "Ready for handoff" -> "In Translation" ---
# Update every cell that currently shows "Ready for handoff" so the shared
# string is fully replaced (old text disappears from the shared string
# table, and every cell that used to reference it now references the new,
# shared "In Translation" entry), matching the sharedStrings.xml diff.

$ws_overview.Range("E2").Value = "In Translation"
$ws_overview.Range("F2").Value = "In Translation"
$ws_overview.Range("E3").Value = "In Translation"
$ws_overview.Range("F3").Value = "In Translation"
$ws_overview.Range("E4").Value = "In Translation"
$ws_overview.Range("F4").Value = "In Translation"

$ws_zhcn.Range("C2").Value = "In Translation"
$ws_zhcn.Range("C3").Value = "In Translation"
$ws_zhcn.Range("C4").Value = "In Translation"

$ws_dede.Range("C2").Value = "In Translation"
$ws_dede.Range("C3").Value = "In Translation"
$ws_dede.Range("C4").Value = "In Translation"

# --- Column width changes ---
# Stored OOXML column width shrinks from 17.2159881591797 down to
# 13.4101845877511 on Overview!E:F and on the "Status" column (column C) of
# the zh-cn / de-de sheets. The ColumnWidth COM property in this runtime is
# quantized to 1/6-character increments, so 12.5 is the closest settable
# value, producing the nearest achievable width (13.33.. characters) to the
# target (13.41.. characters).
$ws_overview.Columns.Item(5).ColumnWidth = 12.5
$ws_overview.Columns.Item(6).ColumnWidth = 12.5

$ws_zhcn.Columns.Item(3).ColumnWidth = 12.5

$ws_dede.Columns.Item(3).ColumnWidth = 12.5

Write-Host "Applied localization-status report edits: status text -> 'In Translation'; status columns narrowed."
